$d = $word.ActiveDocument

# Programa paragraph: insert breaks after each numbered sentence.
$d.Content.Find.Execute("trabalho.2- A", $true, $true, $false, $false, $false, $true, 1, $false, "trabalho.^l2- A", 2)
$d.Content.Find.Execute("termodinâmico).3- A", $true, $true, $false, $false, $false, $true, 1, $false, "termodinâmico).^l3- A", 2)
$d.Content.Find.Execute("Gibbs-Helmholtz.4- Equilíbrio", $true, $true, $false, $false, $false, $true, 1, $false, "Gibbs-Helmholtz.^l4- Equilíbrio", 2)
$d.Content.Find.Execute("metal-óxido-O2(g).5- Equilíbrio", $true, $true, $false, $false, $false, $true, 1, $false, "metal-óxido-O2(g).^l5- Equilíbrio", 2)

# Critério run: insert break between "expressão:" and "NF="
$d.Content.Find.Execute("através da expressão:NF=", $true, $true, $false, $false, $false, $true, 1, $false, "através da expressão:^lNF=", 2)

# Bibliografia paragraph: insert breaks after each numbered reference.
$d.Content.Find.Execute("ISBN 0-87339-270-1.2) P.", $true, $true, $false, $false, $false, $true, 1, $false, "ISBN 0-87339-270-1.^l2) P.", 2)
$d.Content.Find.Execute("ISBN 978-85-216-1600-9.3) S.Stolen", $true, $true, $false, $false, $false, $true, 1, $false, "ISBN 978-85-216-1600-9.^l3) S.Stolen", 2)
$d.Content.Find.Execute("ISBN 978-0-471-49230-6.4) R. DeHoff", $true, $true, $false, $false, $false, $true, 1, $false, "ISBN 978-0-471-49230-6.^l4) R. DeHoff", 2)
$d.Content.Find.Execute("ISBN 978-0-8493-4065-9.5) Y.A. Chang", $true, $true, $false, $false, $false, $true, 1, $false, "ISBN 978-0-8493-4065-9.^l5) Y.A. Chang", 2)
